$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.685.71'
$ws.Range('E2').Value = '  +3.53%  '
$ws.Range('D3').Value = '3.462.09'
$ws.Range('E3').Value = '  +4.20%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.58'
$ws.Range('E5').Value = '  +4.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.69'
$ws.Range('E6').Value = '  +3.79%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.462.88'
$ws.Range('E8').Value = '  +4.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.565'
$ws.Range('E9').Value = '  +6.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.57'
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.125'
$ws.Range('E11').Value = '  +6.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.446'
$ws.Range('E12').Value = '  +2.54%  '
$ws.Range('D13').Value = '4.057.62'
$ws.Range('E13').Value = '  +4.15%  '
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000198'
$ws.Range('E15').Value = '  +9.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.92'
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('D17').Value = '64.695.64'
$ws.Range('E17').Value = '  +3.56%  '
$ws.Range('D18').Value = '3.467.71'
$ws.Range('E18').Value = '  +4.26%  '
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.41'
$ws.Range('E20').Value = '  +4.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '397.98'
$ws.Range('E21').Value = '  +4.01%  '
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('E23').Value = '  +1.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.19'
$ws.Range('E24').Value = '  +3.50%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000122'
$ws.Range('E26').Value = '  +26.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.61'
$ws.Range('E27').Value = '  +8.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.180'
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.04'
$ws.Range('E30').Value = '  +8.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.77'
$ws.Range('E31').Value = '  +5.67%  '
$ws.Range('E32').Value = '  +3.62%  '
$ws.Range('E33').Value = '  +5.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.85'
$ws.Range('E34').Value = '  +3.98%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +4.84%  '
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.84'
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0784'
$ws.Range('E39').Value = '  +7.09%  '
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.66'
$ws.Range('E41').Value = '  +2.81%  '
$ws.Range('D42').Value = '2.911.21'
$ws.Range('E42').Value = '  +2.48%  '
$ws.Range('E43').Value = '  +3.44%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.37'
$ws.Range('E44').Value = '  +4.68%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.44'
$ws.Range('E45').Value = '  +2.53%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.773'
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.82'
$ws.Range('E47').Value = '  +8.09%  '
$ws.Range('E48').Value = '  +5.82%  '
$ws.Range('E49').Value = '  +24.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.58'
$ws.Range('E50').Value = '  +4.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.857'
$ws.Range('E51').Value = '  +5.86%  '
